$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.374.62"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.509.35"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'590.71"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'134.47"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'7.60"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'0.389"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").Value = "4.106.99"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "3.507.44"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "'25.73"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "64.350.48"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'10.00"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").Value = "'394.47"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "'0.578"
$ws.Range("D23").Value = "3.649.52"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'74.69"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").Value = "'2.25"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "'8.26"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -6.61%  "
$ws.Range("E33").Value = "  +6.12%  "
$ws.Range("D34").Value = "3.539.03"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D36").Value = "'23.36"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'5.37"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").Value = "'6.95"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "'167.32"
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").Value = "'0.0787"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.811"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'25.11"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "'6.81"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "2.379.92"
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("D50").Value = "'0.896"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  -0.18%  "
